$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.74'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.06'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.130'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05585'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.497'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.022'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8171'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8484'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1340'
$ws.Range("D10").Style = "Normal"

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06957'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02851'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09397'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001526'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006012'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14OneONE'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006276'
$ws.Range("D16").Style = "Normal"

$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03218'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.743'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04713'
$ws.Range("D23").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001248'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004602'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009600'
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03656'
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = 'BKEXToken'

$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1368'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '40BKEXTokenBKKBestin24h'

$ws.Range("B42").Value = 'KickToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006230'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '41KickTokenKICK'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002489'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007883'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005310'
$ws.Range("D45").Style = "Normal"

$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("D50").Style = "Normal"
